$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New row 10: "Find minimum in rotated sorted array with duplicate elements"
# B10 = S.no. 8
$ws.Range("B10").Value = 8
$ws.Range("B10").HorizontalAlignment = -4131
$ws.Range("B10").VerticalAlignment = -4160

# C10 = "Binary S1 23" (page reference note)
$ws.Range("C10").Value = "Binary S1 23"
$ws.Range("C10").HorizontalAlignment = -4131
$ws.Range("C10").VerticalAlignment = -4160

# D10 = Question text, left/top aligned + wrapped
$ws.Range("D10").Value = "Find minimum in rotated sorted array with duplicate elements"
$ws.Range("D10").HorizontalAlignment = -4131
$ws.Range("D10").VerticalAlignment = -4160
$ws.Range("D10").WrapText = $true

# E10 = Link to the LeetCode problem (hyperlink style, wrap text only)
$leetcodeUrl = "https://leetcode.com/problems/find-minimum-in-rotated-sorted-array-ii/?envType=list&envId=raau48es"
$ws.Hyperlinks.Add($ws.Range("E10"), $leetcodeUrl, [System.Reflection.Missing]::Value, [System.Reflection.Missing]::Value, $leetcodeUrl)
$ws.Range("E10").Value = "Find Minimum in Rotated Sorted Array II - LeetCode"
$ws.Range("E10").HorizontalAlignment = 1
$ws.Range("E10").VerticalAlignment = -4107
$ws.Range("E10").WrapText = $true

# Row height to match the wrapped content of the other data rows
$ws.Rows.Item(10).RowHeight = 43.2

# Update selection to the newly active cell, like the authored workbook
$ws.Range("F10").Select()
